$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.090.04"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.309.88"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D5").Value = "542.55"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "132.02"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "2.308.66"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "23.82"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "2.723.37"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "58.976.02"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "2.312.69"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -3.48%  "
$ws.Range("D21").Value = "312.44"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "62.56"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "0.173"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "1.17"
$ws.Range("E30").Value = "  +6.54%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "169.99"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "0.0₃0738"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").Value = "5.86"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").Value = "0.383"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D36").Value = "17.83"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  +4.34%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "4.06"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "307.53"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "37.66"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "1.51"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "140.92"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").Value = "0.0495"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "18.29"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -0.25%  "
